# Auto-generated edit script to update Chocobo_Profits market-data values
# per the scheduled-runner data refresh described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1349.7042
$ws.Range("J112").Value = 1391.6029
$ws.Range("L112").Value = 4174.8087
$ws.Range("N112").Value = -6390.8087
$ws.Range("H127").Value = 2873.8096
$ws.Range("I127").Value = 840.2857
$ws.Range("J127").Value = 3890.5715
$ws.Range("K127").Value = 2520.8571
$ws.Range("L127").Value = 11671.7145
$ws.Range("M127").Value = 2439.1429
$ws.Range("N127").Value = -21591.7145
$ws.Range("H132").Value = 32387248
$ws.Range("H138").Value = 1768.5385
$ws.Range("I138").Value = 1156.9474
$ws.Range("J138").Value = 3428.5715
$ws.Range("K138").Value = 3470.8422
$ws.Range("L138").Value = 10285.7145
$ws.Range("M138").Value = 1669.1578
$ws.Range("N138").Value = -20565.7145
$ws.Range("H141").Value = 7799
$ws.Range("I141").Value = 11089.182
$ws.Range("J141").Value = 3275
$ws.Range("K141").Value = 33267.546
$ws.Range("L141").Value = 9825
$ws.Range("M141").Value = -28087.546
$ws.Range("N141").Value = -20185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2489.7727
$ws.Range("I132").Value = 1261.0625
$ws.Range("K132").Value = 3783.1875
$ws.Range("M132").Value = -1253.1875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H41").Value = 84900
$ws.Range("J41").Value = 84900
$ws.Range("L41").Value = 84900
$ws.Range("N41").Value = -85676
$ws.Range("H86").Value = 2500.6
$ws.Range("I86").Value = 2357
$ws.Range("J86").Value = 2835.6667
$ws.Range("K86").Value = 2357
$ws.Range("L86").Value = 2835.6667
$ws.Range("M86").Value = -1234
$ws.Range("N86").Value = -5081.6667
$ws.Range("H89").Value = 2500.6
$ws.Range("I89").Value = 2357
$ws.Range("J89").Value = 2835.6667
$ws.Range("K89").Value = 11785
$ws.Range("L89").Value = 14178.3335
$ws.Range("M89").Value = -6169
$ws.Range("N89").Value = -25410.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 187189.31
$ws.Range("I31").Value = 376321.66
$ws.Range("K31").Value = 376321.66
$ws.Range("M31").Value = -376026.66
$ws.Range("H34").Value = 187189.31
$ws.Range("I34").Value = 376321.66
$ws.Range("K34").Value = 376321.66
$ws.Range("M34").Value = -376119.66
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H132").Value = 4066.75
$ws.Range("I132").Value = 3649.3572
$ws.Range("K132").Value = 10948.0716
$ws.Range("M132").Value = -8418.071599999999
$ws.Range("H134").Value = 2078.1875
$ws.Range("I134").Value = 1125.0834
$ws.Range("J134").Value = 4937.5
$ws.Range("K134").Value = 3375.2502
$ws.Range("L134").Value = 14812.5
$ws.Range("M134").Value = -840.2501999999999
$ws.Range("N134").Value = -19882.5
$ws.Range("H135").Value = 31410
$ws.Range("J135").Value = 31410
$ws.Range("L135").Value = 31410
$ws.Range("N135").Value = -41550

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5134.2856
$ws.Range("I56").Value = 5134.2856
$ws.Range("K56").Value = 5134.2856
$ws.Range("M56").Value = -4604.2856
$ws.Range("H68").Value = 3325.7292
$ws.Range("I68").Value = 1223.2142
$ws.Range("J68").Value = 4191.4707
$ws.Range("K68").Value = 3669.6426
$ws.Range("L68").Value = 12574.4121
$ws.Range("M68").Value = -2858.6426
$ws.Range("N68").Value = -14196.4121
$ws.Range("H71").Value = 3325.7292
$ws.Range("I71").Value = 1223.2142
$ws.Range("J71").Value = 4191.4707
$ws.Range("K71").Value = 11008.9278
$ws.Range("L71").Value = 37723.2363
$ws.Range("M71").Value = -6952.927799999999
$ws.Range("N71").Value = -45835.2363
$ws.Range("H107").Value = 9635865
$ws.Range("I107").Value = 351.8
$ws.Range("J107").Value = 18557638
$ws.Range("K107").Value = 1055.4
$ws.Range("L107").Value = 55672914
$ws.Range("M107").Value = 864.5999999999999
$ws.Range("N107").Value = -55676754
$ws.Range("H122").Value = 2377.1304
$ws.Range("I122").Value = 700
$ws.Range("J122").Value = 3455.2856
$ws.Range("K122").Value = 6300
$ws.Range("L122").Value = 31097.5704
$ws.Range("M122").Value = -3850
$ws.Range("N122").Value = -35997.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6012.5195
$ws.Range("I70").Value = 5509.433
$ws.Range("J70").Value = 7788.1177
$ws.Range("K70").Value = 5509.433
$ws.Range("L70").Value = 7788.1177
$ws.Range("M70").Value = -5239.433
$ws.Range("N70").Value = -8328.117699999999
$ws.Range("H73").Value = 6012.5195
$ws.Range("I73").Value = 5509.433
$ws.Range("J73").Value = 7788.1177
$ws.Range("K73").Value = 5509.433
$ws.Range("L73").Value = 7788.1177
$ws.Range("M73").Value = -4573.433
$ws.Range("N73").Value = -9660.117699999999
$ws.Range("H132").Value = 5169.1924
$ws.Range("I132").Value = 4430.3
$ws.Range("J132").Value = 5631
$ws.Range("K132").Value = 13290.9
$ws.Range("L132").Value = 16893
$ws.Range("M132").Value = -10760.9
$ws.Range("N132").Value = -21953

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3127.2727
$ws.Range("I22").Value = 1833.3334
$ws.Range("K22").Value = 1833.3334
$ws.Range("M22").Value = -1538.3334
$ws.Range("H26").Value = 25999.5
$ws.Range("J26").Value = 25999.5
$ws.Range("L26").Value = 25999.5
$ws.Range("N26").Value = -26589.5
$ws.Range("H27").Value = 3127.2727
$ws.Range("I27").Value = 1833.3334
$ws.Range("K27").Value = 1833.3334
$ws.Range("M27").Value = -1726.3334
$ws.Range("H29").Value = 17500
$ws.Range("I29").Value = 10000
$ws.Range("J29").Value = 25000
$ws.Range("K29").Value = 10000
$ws.Range("L29").Value = 25000
$ws.Range("M29").Value = -9705
$ws.Range("N29").Value = -25590
$ws.Range("H31").Value = 3180.0908
$ws.Range("I31").Value = 1128.75
$ws.Range("J31").Value = 8650.333000000001
$ws.Range("K31").Value = 1128.75
$ws.Range("L31").Value = 8650.333000000001
$ws.Range("M31").Value = -880.75
$ws.Range("N31").Value = -9146.333000000001
$ws.Range("H43").Value = 29900
$ws.Range("J43").Value = 29900
$ws.Range("L43").Value = 29900
$ws.Range("N43").Value = -30286
$ws.Range("H68").Value = 1071.0146
$ws.Range("I68").Value = 968.619
$ws.Range("J68").Value = 2361.2
$ws.Range("K68").Value = 968.619
$ws.Range("L68").Value = 2361.2
$ws.Range("M68").Value = -219.619
$ws.Range("N68").Value = -3859.2
$ws.Range("H71").Value = 1071.0146
$ws.Range("I71").Value = 968.619
$ws.Range("J71").Value = 2361.2
$ws.Range("K71").Value = 4843.095
$ws.Range("L71").Value = 11806
$ws.Range("M71").Value = -1099.095
$ws.Range("N71").Value = -19294
$ws.Range("H122").Value = 4926.0625
$ws.Range("I122").Value = 3961.2
$ws.Range("K122").Value = 11883.6
$ws.Range("M122").Value = -9433.599999999999
$ws.Range("H132").Value = 5317
$ws.Range("I132").Value = 4170.769
$ws.Range("K132").Value = 12512.307
$ws.Range("M132").Value = -9982.307000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 20333.334
$ws.Range("I33").Value = 19000
$ws.Range("J33").Value = 21000
$ws.Range("K33").Value = 19000
$ws.Range("L33").Value = 21000
$ws.Range("M33").Value = -18750
$ws.Range("N33").Value = -21500
$ws.Range("H36").Value = 20333.334
$ws.Range("I36").Value = 19000
$ws.Range("J36").Value = 21000
$ws.Range("K36").Value = 19000
$ws.Range("L36").Value = 21000
$ws.Range("M36").Value = -18750
$ws.Range("N36").Value = -21500
$ws.Range("H122").Value = 9666.666999999999
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 9666.666999999999
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 29000.001
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -33900.001
$ws.Range("H132").Value = 14496058
$ws.Range("I132").Value = 2001
$ws.Range("J132").Value = 16670167
$ws.Range("K132").Value = 6003
$ws.Range("L132").Value = 50010501
$ws.Range("M132").Value = -3473
$ws.Range("N132").Value = -50015561

